$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 17.93632866666666
$ws.Range("H2").Value = 53.808986
$ws.Range("I2").Value = 0.1226979812530711
$ws.Range("J2").Value = 0.1347750935001359
$ws.Range("M2").Value = 0.5804443333333333
$ws.Range("N2").Value = 1.741333
$ws.Range("O2").Value = 0.002431273010151717
$ws.Range("P2").Value = 0.002435427107574628
$ws.Range("Q2").Value = 10.41104033537089
$ws.Range("R2").Value = 93.69936301833799
$ws.Range("S2").Value = 0.0002983122902206933
$ws.Range("T2").Value = 0.000328234916136136
$ws.Range("G3").Value = 17.93632866666666
$ws.Range("H3").Value = 53.808986
$ws.Range("I3").Value = 0.1226979812530711
$ws.Range("J3").Value = 0.1347750935001359
$ws.Range("O3").Value = 0.0004752041289926495
$ws.Range("P3").Value = 0.00047601606752829
$ws.Range("Q3").Value = 2.034888444785778
$ws.Range("R3").Value = 18.313996003072
$ws.Range("S3").Value = 0.00005830658731052211
$ws.Range("T3").Value = 0.00006415511000869227
$ws.Range("G4").Value = 17.93632866666666
$ws.Range("H4").Value = 53.808986
$ws.Range("I4").Value = 0.1226979812530711
$ws.Range("J4").Value = 0.1347750935001359
$ws.Range("M4").Value = 136.1000366666667
$ws.Range("N4").Value = 408.30011
$ws.Range("O4").Value = 0.5700742118164518
$ws.Range("P4").Value = 0.5710482463260632
$ws.Range("Q4").Value = 2441.134989198718
$ws.Range("R4").Value = 21970.21490278846
$ws.Range("S4").Value = 0.06994695495431431
$ws.Range("T4").Value = 0.07696308079168379
$ws.Range("G5").Value = 17.93632866666666
$ws.Range("H5").Value = 53.808986
$ws.Range("I5").Value = 0.1226979812530711
$ws.Range("J5").Value = 0.1347750935001359
$ws.Range("M5").Value = 1.221658
$ws.Range("N5").Value = 2.443316
$ws.Range("O5").Value = 0.005117086949542552
$ws.Range("P5").Value = 0.003417220037046797
$ws.Range("Q5").Value = 21.91205940626266
$ws.Range("R5").Value = 131.472356437576
$ws.Range("S5").Value = 0.000627856238605307
$ws.Range("T5").Value = 0.0004605561500035198
$ws.Range("G6").Value = 17.93632866666666
$ws.Range("H6").Value = 53.808986
$ws.Range("I6").Value = 0.1226979812530711
$ws.Range("J6").Value = 0.1347750935001359
$ws.Range("M6").Value = 100.7253213333333
$ws.Range("N6").Value = 302.175964
$ws.Range("O6").Value = 0.4219022240948613
$ws.Range("P6").Value = 0.4226230904617871
$ws.Range("Q6").Value = 1806.642468490278
$ws.Range("R6").Value = 16259.7822164125
$ws.Range("S6").Value = 0.05176655118262031
$ws.Range("T6").Value = 0.05695906653230373
$ws.Range("I7").Value = 0.3929554311523962
$ws.Range("J7").Value = 0.4316338739568692
$ws.Range("M7").Value = 0.5804443333333333
$ws.Range("N7").Value = 1.741333
$ws.Range("O7").Value = 0.002431273010151717
$ws.Range("P7").Value = 0.002435427107574628
$ws.Range("Q7").Value = 33.34264184259556
$ws.Range("R7").Value = 300.08377658336
$ws.Range("S7").Value = 0.0009553819339533523
$ws.Range("T7").Value = 0.00105121283718201
$ws.Range("I8").Value = 0.3929554311523962
$ws.Range("J8").Value = 0.4316338739568692
$ws.Range("O8").Value = 0.0004752041289926495
$ws.Range("P8").Value = 0.00047601606752829
$ws.Range("S8").Value = 0.0001867340433937055
$ws.Range("T8").Value = 0.0002054646592929505
$ws.Range("I9").Value = 0.3929554311523962
$ws.Range("J9").Value = 0.4316338739568692
$ws.Range("M9").Value = 136.1000366666667
$ws.Range("N9").Value = 408.30011
$ws.Range("O9").Value = 0.5700742118164518
$ws.Range("P9").Value = 0.5710482463260632
$ws.Range("Q9").Value = 7818.036143587913
$ws.Range("R9").Value = 70362.32529229121
$ws.Range("S9").Value = 0.2240137576931963
$ws.Range("T9").Value = 0.2464837667779952
$ws.Range("I10").Value = 0.3929554311523962
$ws.Range("J10").Value = 0.4316338739568692
$ws.Range("M10").Value = 1.221658
$ws.Range("N10").Value = 2.443316
$ws.Range("O10").Value = 0.005117086949542552
$ws.Range("P10").Value = 0.003417220037046797
$ws.Range("Q10").Value = 70.17607513578666
$ws.Range("R10").Value = 421.05645081472
$ws.Range("S10").Value = 0.002010787108501793
$ws.Range("T10").Value = 0.001474987922753545
$ws.Range("I11").Value = 0.3929554311523962
$ws.Range("J11").Value = 0.4316338739568692
$ws.Range("M11").Value = 100.7253213333333
$ws.Range("N11").Value = 302.175964
$ws.Range("O11").Value = 0.4219022240948613
$ws.Range("P11").Value = 0.4226230904617871
$ws.Range("Q11").Value = 5785.99552244921
$ws.Range("R11").Value = 52073.95970204289
$ws.Range("S11").Value = 0.1657887703733511
$ws.Range("T11").Value = 0.1824184417596455
$ws.Range("G12").Value = 9.626273333333334
$ws.Range("H12").Value = 28.87882
$ws.Range("I12").Value = 0.06585095126993876
$ws.Range("J12").Value = 0.07233263354328205
$ws.Range("M12").Value = 0.5804443333333333
$ws.Range("N12").Value = 1.741333
$ws.Range("O12").Value = 0.002431273010151717
$ws.Range("P12").Value = 0.002435427107574628
$ws.Range("Q12").Value = 5.587515807451111
$ws.Range("R12").Value = 50.28764226706
$ws.Range("S12").Value = 0.0001601016405154181
$ws.Range("T12").Value = 0.000176160856493571
$ws.Range("G13").Value = 9.626273333333334
$ws.Range("H13").Value = 28.87882
$ws.Range("I13").Value = 0.06585095126993876
$ws.Range("J13").Value = 0.07233263354328205
$ws.Range("O13").Value = 0.0004752041289926495
$ws.Range("P13").Value = 0.00047601606752829
$ws.Range("Q13").Value = 1.092107127182222
$ws.Range("R13").Value = 9.82896414464
$ws.Range("S13").Value = 0.00003129264394156865
$ws.Range("T13").Value = 0.000034431495773238
$ws.Range("G14").Value = 9.626273333333334
$ws.Range("H14").Value = 28.87882
$ws.Range("I14").Value = 0.06585095126993876
$ws.Range("J14").Value = 0.07233263354328205
$ws.Range("M14").Value = 136.1000366666667
$ws.Range("N14").Value = 408.30011
$ws.Range("O14").Value = 0.5700742118164518
$ws.Range("P14").Value = 0.5710482463260632
$ws.Range("Q14").Value = 1310.136153630022
$ws.Range("R14").Value = 11791.2253826702
$ws.Range("S14").Value = 0.03753992914257391
$ws.Range("T14").Value = 0.04130542353703699
$ws.Range("G15").Value = 9.626273333333334
$ws.Range("H15").Value = 28.87882
$ws.Range("I15").Value = 0.06585095126993876
$ws.Range("J15").Value = 0.07233263354328205
$ws.Range("M15").Value = 1.221658
$ws.Range("N15").Value = 2.443316
$ws.Range("O15").Value = 0.005117086949542552
$ws.Range("P15").Value = 0.003417220037046797
$ws.Range("Q15").Value = 11.76001382785333
$ws.Range("R15").Value = 70.56008296712
$ws.Range("S15").Value = 0.0003369650433583662
$ws.Range("T15").Value = 0.0002471765246764667
$ws.Range("G16").Value = 9.626273333333334
$ws.Range("H16").Value = 28.87882
$ws.Range("I16").Value = 0.06585095126993876
$ws.Range("J16").Value = 0.07233263354328205
$ws.Range("M16").Value = 100.7253213333333
$ws.Range("N16").Value = 302.175964
$ws.Range("O16").Value = 0.4219022240948613
$ws.Range("P16").Value = 0.4226230904617871
$ws.Range("Q16").Value = 969.6094747424979
$ws.Range("R16").Value = 8726.485272682481
$ws.Range("S16").Value = 0.02778266279954949
$ws.Range("T16").Value = 0.03056944112930179
$ws.Range("G17").Value = 39.29803649999999
$ws.Range("H17").Value = 78.59607299999999
$ws.Range("I17").Value = 0.2688281328564436
$ws.Range("J17").Value = 0.1968591842135532
$ws.Range("M17").Value = 0.5804443333333333
$ws.Range("N17").Value = 1.741333
$ws.Range("O17").Value = 0.002431273010151717
$ws.Range("P17").Value = 0.002435427107574628
$ws.Range("Q17").Value = 22.8103225975515
$ws.Range("R17").Value = 136.861935585309
$ws.Range("S17").Value = 0.0006535945837833515
$ws.Range("T17").Value = 0.0004794361936087147
$ws.Range("G18").Value = 39.29803649999999
$ws.Range("H18").Value = 78.59607299999999
$ws.Range("I18").Value = 0.2688281328564436
$ws.Range("J18").Value = 0.1968591842135532
$ws.Range("O18").Value = 0.0004752041289926495
$ws.Range("P18").Value = 0.00047601606752829
$ws.Range("Q18").Value = 4.458388439615999
$ws.Range("R18").Value = 26.750330637696
$ws.Range("S18").Value = 0.0001277482387227666
$ws.Range("T18").Value = 0.00009370813472616281
$ws.Range("G19").Value = 39.29803649999999
$ws.Range("H19").Value = 78.59607299999999
$ws.Range("I19").Value = 0.2688281328564436
$ws.Range("J19").Value = 0.1968591842135532
$ws.Range("M19").Value = 136.1000366666667
$ws.Range("N19").Value = 408.30011
$ws.Range("O19").Value = 0.5700742118164518
$ws.Range("P19").Value = 0.5710482463260632
$ws.Range("Q19").Value = 5348.464208578005
$ws.Range("R19").Value = 32090.78525146803
$ws.Range("S19").Value = 0.1532519859522255
$ws.Range("T19").Value = 0.112416091918329
$ws.Range("G20").Value = 39.29803649999999
$ws.Range("H20").Value = 78.59607299999999
$ws.Range("I20").Value = 0.2688281328564436
$ws.Range("J20").Value = 0.1968591842135532
$ws.Range("M20").Value = 1.221658
$ws.Range("N20").Value = 2.443316
$ws.Range("O20").Value = 0.005117086949542552
$ws.Range("P20").Value = 0.003417220037046797
$ws.Range("Q20").Value = 48.00876067451699
$ws.Range("R20").Value = 192.035042698068
$ws.Range("S20").Value = 0.001375616930309599
$ws.Range("T20").Value = 0.0006727111487712404
$ws.Range("G21").Value = 39.29803649999999
$ws.Range("H21").Value = 78.59607299999999
$ws.Range("I21").Value = 0.2688281328564436
$ws.Range("J21").Value = 0.1968591842135532
$ws.Range("M21").Value = 100.7253213333333
$ws.Range("N21").Value = 302.175964
$ws.Range("O21").Value = 0.4219022240948613
$ws.Range("P21").Value = 0.4226230904617871
$ws.Range("Q21").Value = 3958.307354231562
$ws.Range("R21").Value = 23749.84412538937
$ws.Range("S21").Value = 0.1134191871514024
$ws.Range("T21").Value = 0.08319723681811809
$ws.Range("G22").Value = 21.87880766666666
$ws.Range("H22").Value = 65.63642299999999
$ws.Range("I22").Value = 0.1496675034681502
$ws.Range("J22").Value = 0.1643992147861598
$ws.Range("M22").Value = 0.5804443333333333
$ws.Range("N22").Value = 1.741333
$ws.Range("O22").Value = 0.002431273010151717
$ws.Range("P22").Value = 0.002435427107574628
$ws.Range("Q22").Value = 12.69942993020655
$ws.Range("R22").Value = 114.294869371859
$ws.Range("S22").Value = 0.0003638825616789022
$ws.Range("T22").Value = 0.0004003823041541973
$ws.Range("G23").Value = 21.87880766666666
$ws.Range("H23").Value = 65.63642299999999
$ws.Range("I23").Value = 0.1496675034681502
$ws.Range("J23").Value = 0.1643992147861598
$ws.Range("O23").Value = 0.0004752041289926495
$ws.Range("P23").Value = 0.00047601606752829
$ws.Range("Q23").Value = 2.482165315655111
$ws.Range("R23").Value = 22.339487840896
$ws.Range("S23").Value = 0.00007112261562408668
$ws.Range("T23").Value = 0.00007825666772724652
$ws.Range("G24").Value = 21.87880766666666
$ws.Range("H24").Value = 65.63642299999999
$ws.Range("I24").Value = 0.1496675034681502
$ws.Range("J24").Value = 0.1643992147861598
$ws.Range("M24").Value = 136.1000366666667
$ws.Range("N24").Value = 408.30011
$ws.Range("O24").Value = 0.5700742118164518
$ws.Range("P24").Value = 0.5710482463260632
$ws.Range("Q24").Value = 2977.706525656281
$ws.Range("R24").Value = 26799.35873090653
$ws.Range("S24").Value = 0.0853215840741418
$ws.Range("T24").Value = 0.09387988330101839
$ws.Range("G25").Value = 21.87880766666666
$ws.Range("H25").Value = 65.63642299999999
$ws.Range("I25").Value = 0.1496675034681502
$ws.Range("J25").Value = 0.1643992147861598
$ws.Range("M25").Value = 1.221658
$ws.Range("N25").Value = 2.443316
$ws.Range("O25").Value = 0.005117086949542552
$ws.Range("P25").Value = 0.003417220037046797
$ws.Range("Q25").Value = 26.72842041644466
$ws.Range("R25").Value = 160.370522498668
$ws.Range("S25").Value = 0.0007658616287674862
$ws.Range("T25").Value = 0.0005617882908420255
$ws.Range("G26").Value = 21.87880766666666
$ws.Range("H26").Value = 65.63642299999999
$ws.Range("I26").Value = 0.1496675034681502
$ws.Range("J26").Value = 0.1643992147861598
$ws.Range("M26").Value = 100.7253213333333
$ws.Range("N26").Value = 302.175964
$ws.Range("O26").Value = 0.4219022240948613
$ws.Range("P26").Value = 0.4226230904617871
$ws.Range("Q26").Value = 2203.749932615196
$ws.Range("R26").Value = 19833.74939353677
$ws.Range("S26").Value = 0.06314505258793794
$ws.Range("T26").Value = 0.06947890422241799
